$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 48/49 swap: EnergySwap <-> RenderToken (B and C columns) ---
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('B49').Value = 'RenderToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'

# --- E column (Volume) updates: these are safe as plain text (padded with spaces / % sign) ---
$ws.Range('E2').Value = '  -0.30%  '
$ws.Range('E3').Value = '  +0.21%  '
$ws.Range('E4').Value = '  -0.17%  '
$ws.Range('E5').Value = '  -2.44%  '
$ws.Range('E6').Value = '  +0.51%  '
$ws.Range('E7').Value = '  -0.16%  '
$ws.Range('E8').Value = '  -1.50%  '
$ws.Range('E9').Value = '  -3.24%  '
$ws.Range('E10').Value = '  +0.38%  '
$ws.Range('E11').Value = '  +0.48%  '
$ws.Range('E12').Value = '  +3.60%  '
$ws.Range('E13').Value = '  -0.45%  '
$ws.Range('E14').Value = '  +3.07%  '
$ws.Range('E15').Value = '  +0.44%  '
$ws.Range('E16').Value = '  -0.30%  '
$ws.Range('E17').Value = '  +1.79%  '
$ws.Range('E18').Value = '  -0.66%  '
$ws.Range('E19').Value = '  +0.08%  '
$ws.Range('E20').Value = '  -0.30%  '
$ws.Range('E21').Value = '  +19.84%  '
$ws.Range('E22').Value = '  -0.19%  '
$ws.Range('E23').Value = '  -0.01%  '
$ws.Range('E24').Value = '  -0.14%  '
$ws.Range('E25').Value = '  +6.38%  '
$ws.Range('E26').Value = '  +1.90%  '
$ws.Range('E27').Value = '  -0.11%  '
$ws.Range('E28').Value = '  -0.45%  '
$ws.Range('E29').Value = '  -0.53%  '
$ws.Range('E30').Value = '  +2.19%  '
$ws.Range('E31').Value = '  +1.30%  '
$ws.Range('E32').Value = '  +4.21%  '
$ws.Range('E33').Value = '  +1.32%  '
$ws.Range('E34').Value = '  -0.69%  '
$ws.Range('E35').Value = '  -0.46%  '
$ws.Range('E36').Value = '  +1.45%  '
$ws.Range('E37').Value = '  -0.23%  '
$ws.Range('E38').Value = '  -3.18%  '
$ws.Range('E39').Value = '  -0.01%  '
$ws.Range('E40').Value = '  +0.34%  '
$ws.Range('E41').Value = '  +0.08%  '
$ws.Range('E42').Value = '  +1.28%  '
$ws.Range('E43').Value = '  -2.38%  '
$ws.Range('E44').Value = '  -0.14%  '
$ws.Range('E45').Value = '  -0.14%  '
$ws.Range('E46').Value = '  +4.21%  '
$ws.Range('E47').Value = '  +2.23%  '
$ws.Range('E48').Value = '  +3.63%  '
$ws.Range('E49').Value = '  +0.24%  '
$ws.Range('E50').Value = '  -1.57%  '
$ws.Range('E51').Value = '  -0.02%  '

# --- D column (Price) updates that are safe as plain text (contain two dots, not numeric) ---
$ws.Range('D2').Value = '29.817.55'
$ws.Range('D3').Value = '1.894.46'
$ws.Range('D13').Value = '1.883.02'
$ws.Range('D16').Value = '29.834.67'
$ws.Range('D23').Value = '2.145.89'
$ws.Range('D46').Value = '1.032.65'
$ws.Range('D51').Value = '2.045.20'

# --- D column (Price) updates that look numeric: force text without leaving a stray style behind ---
# Use a donor cell's General/default format, paste-special (formats only) after writing text,
# so the cell keeps its original (unstyled) appearance while the content remains a text string.
$donor = $ws.Range('B2')
$donor.Copy()
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.8001'
$ws.Range('D5').PasteSpecial(-4122)
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '242.88'
$ws.Range('D6').PasteSpecial(-4122)
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.000'
$ws.Range('D7').PasteSpecial(-4122)
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3171'
$ws.Range('D8').PasteSpecial(-4122)
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '25.54'
$ws.Range('D9').PasteSpecial(-4122)
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07054'
$ws.Range('D10').PasteSpecial(-4122)
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.08073'
$ws.Range('D11').PasteSpecial(-4122)
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.7718'
$ws.Range('D12').PasteSpecial(-4122)
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.352'
$ws.Range('D14').PasteSpecial(-4122)
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '92.59'
$ws.Range('D15').PasteSpecial(-4122)
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '5.987'
$ws.Range('D17').PasteSpecial(-4122)
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '13.91'
$ws.Range('D18').PasteSpecial(-4122)
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '244.83'
$ws.Range('D19').PasteSpecial(-4122)
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000007727'
$ws.Range('D20').PasteSpecial(-4122)
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '8.275'
$ws.Range('D21').PasteSpecial(-4122)
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.9998'
$ws.Range('D22').PasteSpecial(-4122)
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1658'
$ws.Range('D25').PasteSpecial(-4122)
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.350'
$ws.Range('D26').PasteSpecial(-4122)
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '165.62'
$ws.Range('D27').PasteSpecial(-4122)
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.74'
$ws.Range('D28').PasteSpecial(-4122)
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.059'
$ws.Range('D29').PasteSpecial(-4122)
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.443'
$ws.Range('D32').PasteSpecial(-4122)
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05696'
$ws.Range('D33').PasteSpecial(-4122)
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.040'
$ws.Range('D34').PasteSpecial(-4122)
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7397'
$ws.Range('D36').PasteSpecial(-4122)
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.9987'
$ws.Range('D37').PasteSpecial(-4122)
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.639'
$ws.Range('D38').PasteSpecial(-4122)
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.784'
$ws.Range('D40').PasteSpecial(-4122)
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.4416'
$ws.Range('D41').PasteSpecial(-4122)
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '72.57'
$ws.Range('D42').PasteSpecial(-4122)
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.816'
$ws.Range('D43').PasteSpecial(-4122)
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.0000'
$ws.Range('D45').PasteSpecial(-4122)
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '103.05'
$ws.Range('D47').PasteSpecial(-4122)
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '10.05'
$ws.Range('D48').PasteSpecial(-4122)
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.874'
$ws.Range('D49').PasteSpecial(-4122)
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.434'
$ws.Range('D50').PasteSpecial(-4122)
